$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced back to
# text (the source data stores prices/ids as text, e.g. multi-dot thousands
# separators elsewhere in the same column) - otherwise Excel's normal type
# inference would silently turn "213.07" into the number 213.07.
$textRefs = @("D5", "D8", "D14", "D18", "D20", "D23", "D25", "D29", "D30", "D36", "D38", "D40", "D44", "D45", "D46", "D49", "D51")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.663.82"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.645.13"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "213.07"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "23.02"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "1.875.93"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.640.17"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").Value = "27.635.59"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "229.52"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "10.06"
$ws.Range("E23").Value = "  +7.58%  "
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("D25").Value = "149.15"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "15.64"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").Value = "1.439.29"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "0.882"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").Value = "0.894"
$ws.Range("E40").Value = "  +13.87%  "
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("B44").Value = "mCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "2.26"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "65.34"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "1.785.91"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "86.55"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").Value = "0.0988"
$ws.Range("E51").Value = "  -2.41%  "

# Restore the default (General) style on the cells we had to force to Text,
# so they end up looking exactly like every other unstyled data cell.
foreach ($ref in $textRefs) {
    $ws.Range($ref).Style = "Normal"
}
